$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1: plain (non-shared) formula
$ws.Range("F1").Formula = '=CONCATENATE("insert into [dbo].[parking history] values(",B1,",",E1,",''",D1,"'',''",C1,"'')")'

# F2:F10 use the same relative formula pattern (originally a shared formula)
$ws.Range("F2:F10").Formula = '=CONCATENATE("insert into [dbo].[parking history] values(",B2,",",E2,",''",D2,"'',''",C2,"'')")'

# Update selection to F12
$ws.Range("F12").Select()
